# Vampire Werewolf Game.docx — replace the trailing "31" scratch paragraph
# and the final empty paragraph with the new "Dev:" TODO list.
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count

# Locate the paragraph whose text is just "31" and the paragraph right
# after it (the final, empty paragraph that precedes the section break).
$thirtyOneIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "31") {
        $thirtyOneIndex = $i
        break
    }
}

if ($thirtyOneIndex -eq -1) {
    throw "Could not find the '31' paragraph"
}

$startPara = $d.Paragraphs.Item($thirtyOneIndex)
$endIndex = [Math]::Min($thirtyOneIndex + 1, $count)
$endPara = $d.Paragraphs.Item($endIndex)

$start = $startPara.Range.Start
$end = $endPara.Range.End

$target = $d.Range($start, $end)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Body"/>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
</w:pPr>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Body"/>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
<w:lastRenderedPageBreak/>
<w:t>Dev:</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Body"/>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
<w:tab/>
<w:t>Oyunu 2.5d yapma</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Body"/>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
</w:rPr>
<w:tab/>
<w:t>Can barını kanlı animasyon ve kızaran ekrana çevirme</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
